# Update Name of Algo
# Applies updated numeric results (re-run KNN imputation values) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -8.257
$ws.Range("E3").Value = 16.856
$ws.Range("A12").Value = -21.564
$ws.Range("D14").Value = -7.391000000000001
$ws.Range("E20").Value = 16.306
$ws.Range("E25").Value = 17.102
$ws.Range("D26").Value = -7.689
$ws.Range("A27").Value = -21.84
$ws.Range("E30").Value = 16.214
$ws.Range("D31").Value = -7.915000000000002
$ws.Range("A32").Value = -21.705
$ws.Range("D35").Value = -7.939
$ws.Range("A36").Value = -20.339
$ws.Range("D37").Value = -7.741
$ws.Range("A38").Value = -19.741
$ws.Range("E44").Value = 16.611
$ws.Range("D45").Value = -7.396000000000001
$ws.Range("A46").Value = -21.825
$ws.Range("E47").Value = 16.347
$ws.Range("D52").Value = -7.87
$ws.Range("A54").Value = -21.703
$ws.Range("A55").Value = -22.196
$ws.Range("A56").Value = -22.159
$ws.Range("D57").Value = -8.289999999999999
$ws.Range("E58").Value = 16.65
$ws.Range("A67").Value = -21.6
$ws.Range("A69").Value = -21.68
$ws.Range("A72").Value = -21.445
$ws.Range("E78").Value = 16.391
$ws.Range("D81").Value = -7.322
$ws.Range("A83").Value = -21.795
$ws.Range("D83").Value = -8.289
$ws.Range("E84").Value = 16.384
$ws.Range("A86").Value = -22.035
$ws.Range("E89").Value = 17.289
$ws.Range("A91").Value = -21.661
$ws.Range("E91").Value = 17.243
$ws.Range("E92").Value = 16.891
$ws.Range("A93").Value = -21.288
$ws.Range("E96").Value = 16.424
$ws.Range("A99").Value = -20.43
$ws.Range("D100").Value = -8.238000000000001
$ws.Range("D102").Value = -7.507000000000001
$ws.Range("E102").Value = 16.656
